$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-17 (the 16 term entries) are reordered into a new sequence.
# The underlying term/definition/notes/examples/sources content is unchanged
# per term - only the row position changes. Re-write each row fully (including
# blank columns) so no stale data from the previous occupant of that row remains.

# Row 2: Instrument Data
$ws.Cells.Item(2, 1).Value2 = "Instrument Data"
$ws.Cells.Item(2, 2).Value2 = ""
$ws.Cells.Item(2, 3).Value2 = ""
$ws.Cells.Item(2, 4).Value2 = "Data created by an instrument including scientific measurements and any engineering or ancillary data which may be included in the data packets."
$ws.Cells.Item(2, 5).Value2 = ""
$ws.Cells.Item(2, 6).Value2 = ""
$ws.Cells.Item(2, 7).Value2 = "- [EO Data Stewardship Glossary](https://ceos.org/document_management/Working_Groups/WGISS/Interest_Groups/Data_Stewardship/White_Papers/EO-DataStewardshipGlossary.pdf)"

# Row 3: Geolocating
$ws.Cells.Item(3, 1).Value2 = "Geolocating"
$ws.Cells.Item(3, 2).Value2 = ""
$ws.Cells.Item(3, 3).Value2 = ""
$ws.Cells.Item(3, 4).Value2 = "Determination of the geographic location of a >=2D feature(?)."
$ws.Cells.Item(3, 5).Value2 = ""
$ws.Cells.Item(3, 6).Value2 = ""
$ws.Cells.Item(3, 7).Value2 = "- ISO 19130-1:2018, 3.36 (‘geopositioning’), modified"

# Row 4: Verification
$ws.Cells.Item(4, 1).Value2 = "Verification"
$ws.Cells.Item(4, 2).Value2 = ""
$ws.Cells.Item(4, 3).Value2 = ""
$ws.Cells.Item(4, 4).Value2 = "The evaluation of whether or not a product, service, or system complies with a regulation requirement, specification, or imposed condition. It is often an internal Process."
$ws.Cells.Item(4, 5).Value2 = ""
$ws.Cells.Item(4, 6).Value2 = ""
$ws.Cells.Item(4, 7).Value2 = "- EU-US Land Imaging EO Collaboration"

# Row 5: Auxiliary Data
$ws.Cells.Item(5, 1).Value2 = "Auxiliary Data"
$ws.Cells.Item(5, 2).Value2 = ""
$ws.Cells.Item(5, 3).Value2 = ""
$ws.Cells.Item(5, 4).Value2 = "Data required to perform processing of Sensor Data which is not obtained from the Sensor itself. Include: (a) Data provided by the spacecraft (e.g. orbit Position and velocity, attitude, instrument house-keeping Data, on-board time), (b) Data not available from on-board sources."
$ws.Cells.Item(5, 5).Value2 = "For EnMAP, this includes (a) Orbit files, attitude files, Calibration Data, instrument house-keeping Data, (b) atmospheric parameters, Reference images."
$ws.Cells.Item(5, 6).Value2 = ""
$ws.Cells.Item(5, 7).Value2 = "- ENMAP Glossary of Terms, https://www.enmap.org/Data/doc/EnMAP_Terms.pdf, 20210624`n- EO Data Stewardship Glossary)"

# Row 6: Characteristic
$ws.Cells.Item(6, 1).Value2 = "Characteristic"
$ws.Cells.Item(6, 2).Value2 = "- base"
$ws.Cells.Item(6, 3).Value2 = ""
$ws.Cells.Item(6, 4).Value2 = "Abstraction of a Property of an Object or of a set of objects."
$ws.Cells.Item(6, 5).Value2 = "- Characteristics are used for describing Concepts."
$ws.Cells.Item(6, 6).Value2 = ""
$ws.Cells.Item(6, 7).Value2 = "- ISO 1087-1:2000, 3.2.4; ISO 19146:2010(E); https://www.iso.org/standard/20057.html"

# Row 7: Georectifying
$ws.Cells.Item(7, 1).Value2 = "Georectifying"
$ws.Cells.Item(7, 2).Value2 = ""
$ws.Cells.Item(7, 3).Value2 = "Orthorectifying"
$ws.Cells.Item(7, 4).Value2 = "The correction of sample locations to achieve some sort of geometric regularity, e.g., a regular 2D geographic grid."
$ws.Cells.Item(7, 5).Value2 = ""
$ws.Cells.Item(7, 6).Value2 = ""
$ws.Cells.Item(7, 7).Value2 = "- KCEO"

# Row 8: Test term
$ws.Cells.Item(8, 1).Value2 = "Test term"
$ws.Cells.Item(8, 2).Value2 = ""
$ws.Cells.Item(8, 3).Value2 = ""
$ws.Cells.Item(8, 4).Value2 = "Second_definition_goes_here."
$ws.Cells.Item(8, 5).Value2 = "- here should be bullets`n- like this"
$ws.Cells.Item(8, 6).Value2 = "- this is also bullets`n- like this"
$ws.Cells.Item(8, 7).Value2 = "- KCEO (no link included, so no brackets)`n- [Website](https://en.wikipedia.org/wiki/Thai_script) ( if you have web references, just add the term goes into square [] brackets and the url into () normal brackets`n```````n`n---`n`nReferences: `n`n1. Strobl, P. A., Woolliams, E. R., & Molch, K. (2024). Lost in Translation: The Need for Common Vocabularies and an Interoperable Thesaurus in Earth Observation Sciences. Surveys in Geophysics, 1-29."

# Row 9: Quantity
$ws.Cells.Item(9, 1).Value2 = "Quantity"
$ws.Cells.Item(9, 2).Value2 = "- base"
$ws.Cells.Item(9, 3).Value2 = ""
$ws.Cells.Item(9, 4).Value2 = "Property whose instances can be compared by ratio or only by order."
$ws.Cells.Item(9, 5).Value2 = ""
$ws.Cells.Item(9, 6).Value2 = ""
$ws.Cells.Item(9, 7).Value2 = "- gEOGlos(VIM4 Notes omitted)"

# Row 10: Ancillary Data
$ws.Cells.Item(10, 1).Value2 = "Ancillary Data"
$ws.Cells.Item(10, 2).Value2 = ""
$ws.Cells.Item(10, 3).Value2 = ""
$ws.Cells.Item(10, 4).Value2 = "Data other than instrument measurements, originating in the instrument itself or from the satellite, required to perform processing of the Data. They include orbit Data, attitude Data, time Information, and spacecraft engineering Data, Calibration Data, Data quality Information, and Data from other instruments or earth system models."
$ws.Cells.Item(10, 5).Value2 = ""
$ws.Cells.Item(10, 6).Value2 = ""
$ws.Cells.Item(10, 7).Value2 = "- CEOS-ARD PFS template 20220302"

# Row 11: Validation
$ws.Cells.Item(11, 1).Value2 = "Validation"
$ws.Cells.Item(11, 2).Value2 = ""
$ws.Cells.Item(11, 3).Value2 = ""
$ws.Cells.Item(11, 4).Value2 = "Validation aims to verify that the specified requirements are achieved or compliant. This involves comparing  mission products with representative Reference Data, considering various Observation conditions, ensuring the quality and Traceability of the Reference Data used."
$ws.Cells.Item(11, 5).Value2 = "- In this part of ISO 19159, the term validation is used in a limited sense and only relates to the validation of Calibration Data in order to control their change over time."
$ws.Cells.Item(11, 6).Value2 = ""
$ws.Cells.Item(11, 7).Value2 = "- BIPM; QA4EO; ESA ?, modified"

# Row 12: Entity
$ws.Cells.Item(12, 1).Value2 = "Entity"
$ws.Cells.Item(12, 2).Value2 = "- base"
$ws.Cells.Item(12, 3).Value2 = ""
$ws.Cells.Item(12, 4).Value2 = "A government or business organization that is formed to conduct business or represent the government of the day."
$ws.Cells.Item(12, 5).Value2 = ""
$ws.Cells.Item(12, 6).Value2 = "CEOS Entities include Working Groups, Virtual Constellations, etc."
$ws.Cells.Item(12, 7).Value2 = "- WGISS Shared Collection Lifecycle Management Principles for Earth Observation Data)"

# Row 13: Data
$ws.Cells.Item(13, 1).Value2 = "Data"
$ws.Cells.Item(13, 2).Value2 = "- core"
$ws.Cells.Item(13, 3).Value2 = ""
$ws.Cells.Item(13, 4).Value2 = "Scientific or technical measurements, values calculated therefrom, observations, or facts that can be represented by numbers, tables, graphs, models, text, or symbols which are used as a basis for reasoning and further calculation."
$ws.Cells.Item(13, 5).Value2 = ""
$ws.Cells.Item(13, 6).Value2 = ""
$ws.Cells.Item(13, 7).Value2 = "- WGISS Shared Collection Lifecycle Management Principles for Earth Observation Data)"

# Row 14: Baseline
$ws.Cells.Item(14, 1).Value2 = "Baseline"
$ws.Cells.Item(14, 2).Value2 = ""
$ws.Cells.Item(14, 3).Value2 = ""
$ws.Cells.Item(14, 4).Value2 = "Source data that has been processed to a common set of requirements and organised into a form that allows immediate analysis and interoperability through time and with other collections."
$ws.Cells.Item(14, 5).Value2 = ""
$ws.Cells.Item(14, 6).Value2 = ""
$ws.Cells.Item(14, 7).Value2 = "- WGISS Shared Collection Lifecycle Management Principles for Earth Observation Data)"

# Row 15: Reference
$ws.Cells.Item(15, 1).Value2 = "Reference"
$ws.Cells.Item(15, 2).Value2 = ""
$ws.Cells.Item(15, 3).Value2 = ""
$ws.Cells.Item(15, 4).Value2 = "A sort of Data acquired with an Uncertainty significantly lower (quantify?) than that of the Data it is being compared with."
$ws.Cells.Item(15, 5).Value2 = ""
$ws.Cells.Item(15, 6).Value2 = ""
$ws.Cells.Item(15, 7).Value2 = "- VIM?, modified"

# Row 16: User
$ws.Cells.Item(16, 1).Value2 = "User"
$ws.Cells.Item(16, 2).Value2 = ""
$ws.Cells.Item(16, 3).Value2 = ""
$ws.Cells.Item(16, 4).Value2 = "External person, institution or system that consumes provided services."
$ws.Cells.Item(16, 5).Value2 = "Includes Data Access or Science and Service Exploitation Platforms provided by a payload data ground segment."
$ws.Cells.Item(16, 6).Value2 = ""
$ws.Cells.Item(16, 7).Value2 = "- EO Data Stewardship Glossary)"

# Row 17: Uncertainty
$ws.Cells.Item(17, 1).Value2 = "Uncertainty"
$ws.Cells.Item(17, 2).Value2 = ""
$ws.Cells.Item(17, 3).Value2 = ""
$ws.Cells.Item(17, 4).Value2 = "Non-negative parameter, associated with Data, which characterizes the dispersion of the values of a [Trait ]that could reasonably be attributed to a Phenomenon [by means of sensing or modelling]."
$ws.Cells.Item(17, 5).Value2 = "- In case of quantitative(continuous) Data the uncertainty may be, for example, a standard deviation (or a given multiple of it), or the half-width of an interval having a stated level of confidence. (see e.g. standard and Expanded uncertainty)`n- For qualitative (categorical?) Data uncertainty may be, for example, expressed by commission and omission (‘confusion matrix’) or overall errors."
$ws.Cells.Item(17, 6).Value2 = ""
$ws.Cells.Item(17, 7).Value2 = "- modified from GUM, VIM4 :3.1, FIDUCEO, Notes added"

Write-Host "Reordered 16 glossary rows."
